# Updates the Neo4j/StatQuery text for TC05 Canine StudyUBC02 workbook:
# remove the 'Bladder' entry from each primary_disease_site filter list,
# and drop the Cohort column/line from the Cases query (B2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$casesQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)

MATCH (c)<--(diag:diagnosis)
WHERE s.clinical_study_designation IN ['UBC02'] and demo.breed in ['Mixed Breed', 'Scottish Terrier','Shetland Sheepdog']and diag.disease_term in ['Bladder Cancer','Healthy Control'] and diag.primary_disease_site in ['Bladder, Urethra', 'Bladder, Urethra, Prostate', 'Bladder, Urethra, Vagina', 'Not Applicable']
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`
'@

$statQuery = @'
MATCH (s:study)
  MATCH (demo:demographic) 
  MATCH (diag:diagnosis)
 MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis)
	WHERE s.clinical_study_designation IN ['UBC02'] and demo.breed in ['Mixed Breed', 'Scottish Terrier','Shetland Sheepdog']and diag.disease_term in ['Bladder Cancer','Healthy Control'] and diag.primary_disease_site in ['Bladder, Urethra', 'Bladder, Urethra, Prostate', 'Bladder, Urethra, Vagina', 'Not Applicable']
    
OPTIONAL MATCH (s)<-[:member_of]-(c:case)
OPTIONAL MATCH (c)<-[:of_case]-(samp:sample)<-[:of_sample]-(f:file)
RETURN 
	count(DISTINCT(f)) as number_of_files , 
	count(DISTINCT(samp)) as number_of_sample , 
	count(DISTINCT(c.case_id)) as number_of_cases , 
	count(DISTINCT(s.clinical_study_designation)) as number_of_study
'@

$samplesQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic), (samp:sample)-->(c)<--(diag:diagnosis) 
WHERE s.clinical_study_designation IN ['UBC02'] and demo.breed in ['Mixed Breed', 'Scottish Terrier','Shetland Sheepdog']and diag.disease_term in ['Bladder Cancer','Healthy Control'] and diag.primary_disease_site in ['Bladder, Urethra', 'Bladder, Urethra, Prostate', 'Bladder, Urethra, Vagina', 'Not Applicable']
    

 WITH DISTINCT samp AS samp, c, demo, diag
RETURN  coalesce(samp.sample_id, '') AS `Sample ID`, 
        coalesce(c.case_id, '') AS `Case ID`, 
        coalesce(demo.breed,'') AS Breed , 
        coalesce(diag.disease_term,'') AS Diagnosis , 
        coalesce(samp.sample_site, '') AS `Sample Site`,
        coalesce(samp.summarized_sample_type, '') AS `Sample Type`,
        coalesce(samp.specific_sample_pathology, '') AS `Pathology/Morphology`,
        coalesce(samp.tumor_grade, '') AS `Tumor Grade`,
        coalesce(samp.sample_chronology, '') AS `Sample Chronology`,
        coalesce(samp.percentage_tumor, '') AS `Percentage Tumor`,
        coalesce(samp.necropsy_sample, '') AS `Necropsy Sample`,
        coalesce(samp.sample_preservation, '') AS `Sample Preservation`
'@

$filesQuery = @'
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
 MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
WHERE s.clinical_study_designation IN ['UBC02'] and demo.breed in ['Mixed Breed', 'Scottish Terrier','Shetland Sheepdog']and diag.disease_term in ['Bladder Cancer','Healthy Control'] and diag.primary_disease_site in ['Bladder, Urethra', 'Bladder, Urethra, Prostate', 'Bladder, Urethra, Vagina', 'Not Applicable']
    
WITH DISTINCT f, parent, c, demo, diag, s
RETURN coalesce(f.file_name, '') AS `File Name`, 
        coalesce(f.file_type, '') AS `File Type`, 
        coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(f.file_format, '') AS `File Format`,
        coalesce(f.file_size, '') AS `Size`,
        coalesce(c.case_id, '') AS `Case ID`, 
        coalesce(demo.breed,'') AS Breed , 
        coalesce(diag.disease_term,'') AS Diagnosis , 
        coalesce(s.clinical_study_designation,'') AS `Study Code`
'@

# Row 2 - CasesTab
$ws.Range("B2").Value = $casesQuery
$ws.Range("C2").Value = $statQuery

# Row 3 - SamplesTab
$ws.Range("B3").Value = $samplesQuery
$ws.Range("C3").Value = $statQuery

# Row 4 - FilesTab
$ws.Range("B4").Value = $filesQuery
$ws.Range("C4").Value = $statQuery

# Row heights follow Excel's wrap-text autofit for the shorter query text
$ws.Rows.Item(2).RowHeight = 304.5
$ws.Rows.Item(3).RowHeight = 304.5
$ws.Rows.Item(4).RowHeight = 275.5

# Match the author's final selection/scroll position
$ws.Range("C4").Select()
